# df_over_under_summary_forslide.xlsx edit
# "show the EqualOver specific bid in Winter_16 OFF."
#
# The "No New Path" table (rows 3-5) is updated to the same bid figures as
# the "With New Path" table (rows 9-11), and both tables are refreshed with
# the new EqualOver-specific numbers. A new comparison block (cols G:K,
# rows 4-6) is added next to the "No New Path" table showing the spread
# between row 4/row 3 and row 5/row 3 (dollar amounts, bid counts, and the
# ratio between the two).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- "No New Path" table (rows 3-5) ----
$ws.Range("B3").Value = -80512745.606122807
$ws.Range("C3").Value = 29065515.726577502
$ws.Range("D3").Value = 16366
$ws.Range("E3").Value = 4545

$ws.Range("B4").Value = -35545412.737436697
$ws.Range("C4").Value = 38112052.754308701
$ws.Range("D4").Value = 12684
$ws.Range("E4").Value = 5918

$ws.Range("B5").Value = -41092204.375834197
$ws.Range("C5").Value = 37258708.666714497
$ws.Range("D5").Value = 13279
$ws.Range("E5").Value = 5655

# ---- "With New Path" table (rows 9-11) ----
$ws.Range("B9").Value = -80512745.606122807
$ws.Range("C9").Value = 29065515.726577502
$ws.Range("D9").Value = 16366
$ws.Range("E9").Value = 4545

$ws.Range("B10").Value = -35545412.737436697
$ws.Range("C10").Value = 38112052.754308701
$ws.Range("D10").Value = 12684
$ws.Range("E10").Value = 5918

$ws.Range("B11").Value = -41092204.375834197
$ws.Range("C11").Value = 37258708.666714497
$ws.Range("D11").Value = 13279
$ws.Range("E11").Value = 5655

# ---- New comparison block next to the "No New Path" table ----
$ws.Range("G4").Formula = "=B4-B`$3"
$ws.Range("H4").Formula = "=C4-C`$3"
$ws.Range("I4").Formula = "=D4-D`$3"
$ws.Range("J4").Formula = "=E4-E`$3"
$ws.Range("K4").Formula = "=G4/H4"

$ws.Range("G5").Formula = "=B5-B`$3"
$ws.Range("H5").Formula = "=C5-C`$3"
$ws.Range("I5").Formula = "=D5-D`$3"
$ws.Range("J5").Formula = "=E5-E`$3"
$ws.Range("K5").Formula = "=G5/H5"

$ws.Range("G6").Formula = "=G4/G5"

# Match number formats of the new cells to their neighboring columns
# (same style as columns B/C and D/E respectively).
$ws.Range("G4:H5").NumberFormat = '_("$"* #,##0_);_("$"* \(#,##0\);_("$"* "-"??_);_(@_)'
$ws.Range("I4:J5").NumberFormat = '_(* #,##0_);_(* \(#,##0\);_(* "-"??_);_(@_)'

# Column widths for the newly used columns.
$ws.Columns.Item(7).ColumnWidth = 12.5703125
$ws.Columns.Item(8).ColumnWidth = 15.140625
$ws.Columns.Item(11).ColumnWidth = 12

# Selection as left by the author.
$ws.Range("E8").Select() | Out-Null
